# Append the new "Clarida, Gali, Gertler 1999" record as row 9 of Sheet1,
# then move the selection on to C10 (one row below the new data), matching
# the author's editing session captured in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Clarida, Gali, Gertler 1999"
$ws.Range("B9").Value = "p. 1701"
$ws.Range("C9").Value = "`"… if the rule calls for an overly aggressive response of interest rates to movements in expected inflation. In this instance, there is a 'policy overkill' effect that emerges and may result in an oscillating equilibrium.`""

$ws.Range("C10").Select()
